# Insert a new data row at row 283 (pushing the existing rows 283:394 down
# to 284:395) and populate the new row with the new price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 283; this shifts rows 283:394
# down to 284:395 and copies formatting from the row above (row 282).
$ws.Rows.Item(283).Insert()

# Populate the newly inserted row 283 with the new record's data.
$ws.Cells.Item(283, 1).Value = 8
$ws.Cells.Item(283, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(283, 3).Value = "Coquimbo"
$ws.Cells.Item(283, 4).Value = 45009
$ws.Cells.Item(283, 5).Value = 4
$ws.Cells.Item(283, 6).Value = 100112012
$ws.Cells.Item(283, 7).Value = "Espinaca"
$ws.Cells.Item(283, 8).Value = "Sin especificar"
$ws.Cells.Item(283, 9).Value = "Primera"
$ws.Cells.Item(283, 10).Value = 1960
$ws.Cells.Item(283, 11).Value = 450
$ws.Cells.Item(283, 12).Value = 500
$ws.Cells.Item(283, 13).Value = 475
$ws.Cells.Item(283, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(283, 15).Value = "Provincia del Elqu" + [char]0x00ED
$ws.Cells.Item(283, 16).Value = 950
$ws.Cells.Item(283, 17).Value = 0.5
$ws.Cells.Item(283, 18).Value = "Hortaliza"
